$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16
$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# "25" looks numeric, so force it to be stored as text (matching the
# sibling rows, which hold the Value column as an inline/shared string)
# without leaving the cell's style pointing at a new text-number-format.
$valueCell = $ws.Cells.Item($row, 3)
$valueCell.NumberFormat = "@"
$valueCell.Value = "25"
$valueCell.Style = "Normal"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
